# Add I0 (I) and IF (J) columns to the sheet, matching the H-column header style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (border/bold/alignment) from H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J30
$values = @(
    @(1,5),
    @(1,6),
    @(1,5),
    @(1,5),
    @(1,7),
    @(1,6),
    @(1,4),
    @(1,5),
    @(1,6),
    @(1,5),
    @(1,3),
    @(4,5),
    @(6,6),
    @(5,7),
    @(5,6),
    @(7,8),
    @(5,7),
    @(4,5),
    @(4,6),
    @(5,7),
    @(5,7),
    @(5,7),
    @(6,7),
    @(2,5),
    @(3,5),
    @(1,2),
    @(5,6),
    @(1,1),
    @(1,2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
